$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hide rows 2-18 (previously visible detail rows are collapsed)
$ws.Range("A2:A18").EntireRow.Hidden = $true

# New test-data block (rows 21-65).
# Cells are written in the exact order the strings were first authored so the
# shared-string table is rebuilt in the same sequence as the source workbook.
$ws.Cells.Item(21,1).Value = "Экран Авторизация:"
$ws.Cells.Item(22,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(23,1).Value = "Поля экрана доступны для ввода данных, кнопки кликабельны"
$ws.Cells.Item(24,1).Value = "Экран Главного меню:"
$ws.Cells.Item(26,1).Value = "Топлайн клакабелен, кнопки кликабельны"
$ws.Cells.Item(30,1).Value = "Экран Фильтр новостей:"
$ws.Cells.Item(33,1).Value = "Экран Контрольная панель новостей:"
$ws.Cells.Item(41,1).Value = "Ссылки на экране кликабельны и открываются, кнопки кликабельны"
$ws.Cells.Item(45,1).Value = "Экран Фильтрация Жалоб:"
$ws.Cells.Item(42,1).Value = "Экран Жалобы:"
$ws.Cells.Item(39,1).Value = "Экран О нас:"
$ws.Cells.Item(36,1).Value = "Экран Добавление новости:"
$ws.Cells.Item(27,1).Value = "Экран Новости:"
$ws.Cells.Item(47,1).Value = "Кнопки кликабельны, чек-боксы работают"
$ws.Cells.Item(48,1).Value = "Экран Добавление жалобы:"
$ws.Cells.Item(51,1).Value = "Экран Редактирование жалобы:"
$ws.Cells.Item(53,1).Value = "Кнопки кликабельны"
$ws.Cells.Item(57,1).Value = "Экран Цитаты:"
$ws.Cells.Item(60,1).Value = "Экран Добавление комментария к жалобе:"
$ws.Cells.Item(54,1).Value = "Экран Редактирования комментария к жалобе:"
$ws.Cells.Item(63,1).Value = "Экран Редактирование  новости:"

# Remaining cells reuse already-registered strings (order is not significant)
$ws.Cells.Item(25,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(28,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(29,1).Value = "Топлайн клакабелен, кнопки кликабельны"
$ws.Cells.Item(31,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(32,1).Value = "Поля экрана доступны для ввода данных, кнопки кликабельны"
$ws.Cells.Item(34,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(35,1).Value = "Топлайн клакабелен, кнопки кликабельны"
$ws.Cells.Item(37,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(38,1).Value = "Поля экрана доступны для ввода данных, кнопки кликабельны"
$ws.Cells.Item(40,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(43,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(44,1).Value = "Топлайн клакабелен, кнопки кликабельны"
$ws.Cells.Item(46,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(49,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(50,1).Value = "Поля экрана доступны для ввода данных, кнопки кликабельны"
$ws.Cells.Item(52,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(55,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(56,1).Value = "Поля экрана доступны для ввода данных, кнопки кликабельны"
$ws.Cells.Item(58,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(59,1).Value = "Кнопки кликабельны"
$ws.Cells.Item(61,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(62,1).Value = "Поля экрана доступны для ввода данных, кнопки кликабельны"
$ws.Cells.Item(64,1).Value = "Элементы экрана отображаются"
$ws.Cells.Item(65,1).Value = "Поля экрана доступны для ввода данных, кнопки кликабельны"

# Apply the standard "blank result" style (thin border + green fill) to column B cells
foreach ($r in @(22,23,25,26,28,29,31,32,34,35,37,38,40,41,43,44,46,47,49,50,52,53,55,56,58,59,61,62,64,65)) {
    $c = $ws.Cells.Item($r,2)
    $c.Borders.LineStyle = 1
    $c.Interior.Color = 5287936
}

# Update selected cell to match the saved workbook state
$ws.Range("G39").Select() | Out-Null

Write-Output "done"
